# Auto-generated Excel COM-interop script
# Applies scraped-price updates to the Leve profit tables across sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2061.38
$ws.Range("I32").Value = 1887.2391
$ws.Range("J32").Value = 4064
$ws.Range("K32").Value = 1887.2391
$ws.Range("L32").Value = 4064
$ws.Range("M32").Value = -1600.2391
$ws.Range("N32").Value = -4638

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 710.3095
$ws.Range("I61").Value = 710.3095
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 710.3095
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -498.3095
$ws.Range("N61").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2175.3865
$ws.Range("I132").Value = 2354.8276
$ws.Range("J132").Value = 1828.4667
$ws.Range("K132").Value = 7064.4828
$ws.Range("L132").Value = 5485.4001
$ws.Range("M132").Value = -4534.4828
$ws.Range("N132").Value = -10545.4001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 710.3095
$ws.Range("I136").Value = 710.3095
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 2130.9285
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = 419.0715
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 27742.4
$ws.Range("J2").Value = 27742.4
$ws.Range("L2").Value = 27742.4
$ws.Range("N2").Value = -27968.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 30203.355
$ws.Range("I31").Value = 37166.137
$ws.Range("J31").Value = 17583.312
$ws.Range("K31").Value = 37166.137
$ws.Range("L31").Value = 17583.312
$ws.Range("M31").Value = -36871.137
$ws.Range("N31").Value = -18173.312

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 30203.355
$ws.Range("I34").Value = 37166.137
$ws.Range("J34").Value = 17583.312
$ws.Range("K34").Value = 37166.137
$ws.Range("L34").Value = 17583.312
$ws.Range("M34").Value = -36964.137
$ws.Range("N34").Value = -17987.312

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1058.8206
$ws.Range("I132").Value = 943.9545000000001
$ws.Range("J132").Value = 1207.4706
$ws.Range("K132").Value = 2831.8635
$ws.Range("L132").Value = 3622.4118
$ws.Range("M132").Value = -301.8635000000004
$ws.Range("N132").Value = -8682.4118

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 725
$ws.Range("I17").Value = 725
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 2175
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -2006
$ws.Range("N17").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 3148
$ws.Range("J34").Value = 4000
$ws.Range("L34").Value = 12000
$ws.Range("N34").Value = -12168

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 2214.3333
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 2214.3333
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 6642.999899999999
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = -7230.999899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 3090.9092
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 3090.9092
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 9272.7276
$ws.Range("N55").Value = -9626.7276

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 1308.6666
$ws.Range("I97").Value = 242.66667
$ws.Range("J97").Value = 1460.9524
$ws.Range("K97").Value = 728.00001
$ws.Range("L97").Value = 4382.857199999999
$ws.Range("M97").Value = -232.00001
$ws.Range("N97").Value = -5374.857199999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H101").Value = 8643
$ws.Range("J101").Value = 8643
$ws.Range("L101").Value = 25929
$ws.Range("N101").Value = -30797

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H110").Value = 4211
$ws.Range("J110").Value = 5154
$ws.Range("L110").Value = 15462
$ws.Range("N110").Value = -23642

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 813.98553
$ws.Range("I131").Value = 400
$ws.Range("J131").Value = 892.5
$ws.Range("K131").Value = 1200
$ws.Range("L131").Value = 2677.5
$ws.Range("M131").Value = 3840
$ws.Range("N131").Value = -12757.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1916.05
$ws.Range("I122").Value = 1911.6316
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 5734.8948
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -3284.8948
$ws.Range("N122").Value = -10900

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1693.6571
$ws.Range("I132").Value = 1621.871
$ws.Range("J132").Value = 2250
$ws.Range("K132").Value = 4865.613
$ws.Range("L132").Value = 6750
$ws.Range("M132").Value = -2335.613
$ws.Range("N132").Value = -11810

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 553.1429000000001
$ws.Range("I22").Value = 553.1429000000001
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 553.1429000000001
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -258.1429000000001
$ws.Range("N22").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 553.1429000000001
$ws.Range("I27").Value = 553.1429000000001
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 553.1429000000001
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -446.1429000000001
$ws.Range("N27").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5206.8887
$ws.Range("I40").Value = 5557.7144
$ws.Range("J40").Value = 4983.636
$ws.Range("K40").Value = 5557.7144
$ws.Range("L40").Value = 4983.636
$ws.Range("M40").Value = -5421.7144
$ws.Range("N40").Value = -5255.636

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 665.6
$ws.Range("I46").Value = 616
$ws.Range("J46").Value = 740
$ws.Range("K46").Value = 616
$ws.Range("L46").Value = 740
$ws.Range("M46").Value = -428
$ws.Range("N46").Value = -1116

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1880.7333
$ws.Range("I132").Value = 1466.7
$ws.Range("J132").Value = 3950.9
$ws.Range("K132").Value = 4400.1
$ws.Range("L132").Value = 11852.7
$ws.Range("M132").Value = -1870.1
$ws.Range("N132").Value = -16912.7

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2431.25
$ws.Range("I136").Value = 1833.8889
$ws.Range("J136").Value = 3199.2856
$ws.Range("K136").Value = 5501.6667
$ws.Range("L136").Value = 9597.856800000001
$ws.Range("M136").Value = -2951.6667
$ws.Range("N136").Value = -14697.8568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 629.63635
$ws.Range("I132").Value = 472.65958
$ws.Range("J132").Value = 1551.875
$ws.Range("K132").Value = 1417.97874
$ws.Range("L132").Value = 4655.625
$ws.Range("M132").Value = 1112.02126
$ws.Range("N132").Value = -9715.625
